$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the newly-entered duty-roster names (B24:B32) ---
$ws.Range("B24").Value = "山口玲"
$ws.Range("B25").Value = "なし"
$ws.Range("B26").Value = "なし"
$ws.Range("B27").Value = "志塚惇希"
$ws.Range("B28").Value = "なし"
$ws.Range("B29").Value = "なし"
$ws.Range("B30").Value = "なし"
$ws.Range("B31").Value = "山口洸翔"
$ws.Range("B32").Value = "なし"

# --- Formatting tweaks that came along with the edit ---
# B27 picks up the "Arial 10pt black" look already used higher up the list (e.g. B4)
$ws.Range("B4").Copy()
$ws.Range("B27").PasteSpecial(-4122)   # xlPasteFormats

# B28 loses its explicit style (reverts to the sheet default formatting)
$ws.Range("B28").ClearFormats()

$excel.CutCopyMode = $false

# --- View state: the author zoomed in and left the selection on B32 ---
$excel.ActiveWindow.Zoom = 118
$ws.Range("B32").Select()
